# Applies the "Replace-pattern-text-with-normal-text" template edit:
#   1. In the title heading, collapse the "Adventure Works " + "{" + "Cycle" + "}"
#      runs into a single plain-text run "Adventure Works Cycle".
#   2. In the body paragraph, leave the "{Cycle}" pattern markers alone but wrap
#      the (now grammar-flagged) phrases "bi{" and "is located in" with
#      w:proofErr gramStart/gramEnd markers, splitting the surrounding runs
#      accordingly.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Heading paragraph: "Adventure Works {Cycle}" -> "Adventure Works Cycle"
# ---------------------------------------------------------------------------
$heading = $d.Paragraphs(1).Range
$heading.Find.Execute("Adventure Works {Cycle}", $false, $false, $false, $false, $false, $true, 1, $false, "Adventure Works Cycle", 2)

# ---------------------------------------------------------------------------
# 2) Body paragraph: insert proofErr (gramStart/gramEnd) pairs around
#    "bi{" and "is located in", splitting the existing runs at those points.
# ---------------------------------------------------------------------------
$bodyPara = $d.Paragraphs(2).Range
$bodyPara.MoveEnd(1, -1)
$bodyPara.Text = ""
$bodyPara.Collapse(1)

$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'

$xml = '<w:p ' + $ns + ' w14:paraId="6A4782A5" w14:textId="0FFC0400" w:rsidR="00310EBF" w:rsidRDefault="007A19D5">' +
  '<w:pPr><w:ind w:firstLine="720"/><w:jc w:val="both"/></w:pPr>' +
  '<w:r><w:t xml:space="preserve">Adventure Works </w:t></w:r>' +
  '<w:r w:rsidR="00F14FA6"><w:t>{</w:t></w:r>' +
  '<w:r w:rsidR="008B5ACB"><w:t>Cycle</w:t></w:r>' +
  '<w:r w:rsidR="00F14FA6"><w:t>}</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">, the fictitious company on which the AdventureWorks sample databases are based, is a large, multinational manufacturing company. The company manufactures and sells metal and composite </w:t></w:r>' +
  '<w:proofErr w:type="gramStart"/>' +
  '<w:r><w:t>bi</w:t></w:r>' +
  '<w:r w:rsidR="00F14FA6"><w:t>{</w:t></w:r>' +
  '<w:proofErr w:type="gramEnd"/>' +
  '<w:r w:rsidR="008B5ACB"><w:t>Cycle</w:t></w:r>' +
  '<w:r w:rsidR="00F14FA6"><w:t>}</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> to North American, European and Asian commercial markets. While its base operation </w:t></w:r>' +
  '<w:proofErr w:type="gramStart"/>' +
  '<w:r><w:t>is located in</w:t></w:r>' +
  '<w:proofErr w:type="gramEnd"/>' +
  '<w:r><w:t xml:space="preserve"> Bothell, Washington with 290 employees, several regional sales teams are located throughout their market base.</w:t></w:r>' +
  '</w:p>'

$bodyPara.InsertXML($xml)
